# Person-Hours Estimate.docx edit
#
# Before: a single paragraph containing only the run "120".
# After : that paragraph becomes three runs
#           "Person Hours Estimated for Project 3: " / "120" / " hours"
#         followed by a blank paragraph and a new write-up paragraph.
#
# Plain InsertBefore/InsertAfter text mutations get coalesced back into a
# single <w:r> on save, so we build the exact run/paragraph structure with
# Range.InsertXML (WordprocessingML "single XML" package form), which
# replaces the target range's contents with literally the OOXML supplied.

$d = $word.ActiveDocument

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$bodyFragment = @'
<w:p><w:r><w:t xml:space="preserve">Person Hours Estimated for Project 3: </w:t></w:r><w:r><w:t>120</w:t></w:r><w:r><w:t xml:space="preserve"> hours</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>We laid out the class diagram before coming up with the estimate on this project. Based on the complexity of our class diagram and our experience from the previous two projects, we arrived at the estimate of 120 person hours.</w:t></w:r></w:p>
'@

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document ' + $ns + '><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData>' +
  '</pkg:part></pkg:package>'

# Target the lone paragraph's text ("120"), not its trailing paragraph mark,
# so InsertXML replaces just that run and leaves the rest of the body (the
# sectPr) alone.
$target = $d.Paragraphs(1).Range
$target.MoveEnd(1, -1) | Out-Null
$target.InsertXML($packageXml) | Out-Null
